$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append account 50808 to the client list for the "PREMIUM" group (row 3, column B)
$current = $ws.Range("B3").Value2
$ws.Range("B3").Value = $current + ".50808"

# Update the active selection shown in the sheet view from B7 to B4
$ws.Range("B4").Select()
